$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112043804
$ws.Range("B2").Value = 56398
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 547645
$ws.Range("R2").Value = 6960223
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").Value = "ringhack"

# Row 3
$ws.Range("A3").Value = 112043814
$ws.Range("B3").Value = 90678
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4366
$ws.Range("F3").Value = "Skarp dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum peckii"
$ws.Range("H3").Value = "Banker"
$ws.Range("Q3").Value = 547485
$ws.Range("R3").Value = 6960144
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").ClearContents()

# Row 4
$ws.Range("A4").Value = 112043824
$ws.Range("B4").Value = 78578
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 547465
$ws.Range("R4").Value = 6960199
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").ClearContents()

# Row 5
$ws.Range("A5").Value = 112043834
$ws.Range("B5").Value = 78578
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("Q5").Value = 547645
$ws.Range("R5").Value = 6960007
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").ClearContents()

# Row 6
$ws.Range("A6").Value = 112043803
$ws.Range("B6").Value = 89405
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 547601
$ws.Range("R6").Value = 6959986
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").ClearContents()

# Row 7
$ws.Range("A7").Value = 112043821
$ws.Range("B7").Value = 78578
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6458
$ws.Range("F7").Value = "Lunglav"
$ws.Range("G7").Value = "Lobaria pulmonaria"
$ws.Range("H7").Value = "(L.) Hoffm."
$ws.Range("Q7").Value = 548029
$ws.Range("R7").Value = 6960148
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").ClearContents()

# Row 8
$ws.Range("A8").Value = 112043830
$ws.Range("B8").Value = 78578
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6458
$ws.Range("F8").Value = "Lunglav"
$ws.Range("G8").Value = "Lobaria pulmonaria"
$ws.Range("H8").Value = "(L.) Hoffm."
$ws.Range("Q8").Value = 547485
$ws.Range("R8").Value = 6960130
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 112043838
$ws.Range("B9").Value = 78578
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 548089
$ws.Range("R9").Value = 6960174
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
$ws.Range("AC9").ClearContents()

# Row 10
$ws.Range("A10").Value = 112043827
$ws.Range("B10").Value = 78578
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6458
$ws.Range("F10").Value = "Lunglav"
$ws.Range("G10").Value = "Lobaria pulmonaria"
$ws.Range("H10").Value = "(L.) Hoffm."
$ws.Range("Q10").Value = 547486
$ws.Range("R10").Value = 6960156
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()
$ws.Range("AC10").ClearContents()

# Row 11
$ws.Range("A11").Value = 112043822
$ws.Range("B11").Value = 78578
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6458
$ws.Range("F11").Value = "Lunglav"
$ws.Range("G11").Value = "Lobaria pulmonaria"
$ws.Range("H11").Value = "(L.) Hoffm."
$ws.Range("Q11").Value = 547717
$ws.Range("R11").Value = 6960058
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 112043850
$ws.Range("B12").Value = 89965
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 760
$ws.Range("F12").Value = "Doftticka"
$ws.Range("G12").Value = "Haploporus odorus"
$ws.Range("H12").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q12").Value = 547676
$ws.Range("R12").Value = 6960327
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()
$ws.Range("AC12").ClearContents()

# Row 13
$ws.Range("A13").Value = 112043826
$ws.Range("B13").Value = 78578
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6458
$ws.Range("F13").Value = "Lunglav"
$ws.Range("G13").Value = "Lobaria pulmonaria"
$ws.Range("H13").Value = "(L.) Hoffm."
$ws.Range("Q13").Value = 547407
$ws.Range("R13").Value = 6960191
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()
$ws.Range("AC13").ClearContents()

# Row 14
$ws.Range("A14").Value = 112043831
$ws.Range("B14").Value = 78578
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6458
$ws.Range("F14").Value = "Lunglav"
$ws.Range("G14").Value = "Lobaria pulmonaria"
$ws.Range("H14").Value = "(L.) Hoffm."
$ws.Range("Q14").Value = 547495
$ws.Range("R14").Value = 6960111
$ws.Range("Z14").ClearContents()
$ws.Range("AB14").ClearContents()
$ws.Range("AC14").ClearContents()

# Row 15
$ws.Range("A15").Value = 112043805
$ws.Range("B15").Value = 56398
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = "Tretåig hackspett"
$ws.Range("G15").Value = "Picoides tridactylus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("Q15").Value = 547818
$ws.Range("R15").Value = 6960196
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()
$ws.Range("AC15").Value = "ringhack"

# Row 16
$ws.Range("A16").Value = 112043809
$ws.Range("B16").Value = 89425
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5442
$ws.Range("F16").Value = "Tallticka"
$ws.Range("G16").Value = "Porodaedalea pini"
$ws.Range("H16").Value = "(Brot.) Murrill"
$ws.Range("Q16").Value = 547660
$ws.Range("R16").Value = 6960220
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()
$ws.Range("AC16").ClearContents()

# Row 17
$ws.Range("A17").Value = 112043813
$ws.Range("B17").Value = 78579
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 2081
$ws.Range("F17").Value = "Skrovellav"
$ws.Range("G17").Value = "Lobaria scrobiculata"
$ws.Range("H17").Value = "(Scop.) DC."
$ws.Range("Q17").Value = 548087
$ws.Range("R17").Value = 6960167
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()
$ws.Range("AC17").ClearContents()

# Row 18
$ws.Range("A18").Value = 112043820
$ws.Range("B18").Value = 78578
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6458
$ws.Range("F18").Value = "Lunglav"
$ws.Range("G18").Value = "Lobaria pulmonaria"
$ws.Range("H18").Value = "(L.) Hoffm."
$ws.Range("Q18").Value = 548039
$ws.Range("R18").Value = 6960182
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
$ws.Range("AC18").ClearContents()

# Row 19
$ws.Range("A19").Value = 112043833
$ws.Range("B19").Value = 78578
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6458
$ws.Range("F19").Value = "Lunglav"
$ws.Range("G19").Value = "Lobaria pulmonaria"
$ws.Range("H19").Value = "(L.) Hoffm."
$ws.Range("Q19").Value = 547587
$ws.Range("R19").Value = 6959961
$ws.Range("Z19").ClearContents()
$ws.Range("AB19").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("AC19").ClearContents()

# Row 20
$ws.Range("A20").Value = 112043835
$ws.Range("B20").Value = 78578
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6458
$ws.Range("F20").Value = "Lunglav"
$ws.Range("G20").Value = "Lobaria pulmonaria"
$ws.Range("H20").Value = "(L.) Hoffm."
$ws.Range("Q20").Value = 547724
$ws.Range("R20").Value = 6960020
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()
$ws.Range("AC20").ClearContents()

# Row 21
$ws.Range("A21").Value = 112043817
$ws.Range("B21").Value = 78578
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 6458
$ws.Range("F21").Value = "Lunglav"
$ws.Range("G21").Value = "Lobaria pulmonaria"
$ws.Range("H21").Value = "(L.) Hoffm."
$ws.Range("Q21").Value = 547495
$ws.Range("R21").Value = 6960255
$ws.Range("Z21").ClearContents()
$ws.Range("AB21").ClearContents()
$ws.Range("AC21").ClearContents()

# Row 22
$ws.Range("A22").Value = 112043810
$ws.Range("B22").Value = 89425
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 5442
$ws.Range("F22").Value = "Tallticka"
$ws.Range("G22").Value = "Porodaedalea pini"
$ws.Range("H22").Value = "(Brot.) Murrill"
$ws.Range("Q22").Value = 547840
$ws.Range("R22").Value = 6960201
$ws.Range("Z22").ClearContents()
$ws.Range("AB22").ClearContents()
$ws.Range("AC22").ClearContents()

# Row 23
$ws.Range("A23").Value = 112043837
$ws.Range("B23").Value = 78578
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 6458
$ws.Range("F23").Value = "Lunglav"
$ws.Range("G23").Value = "Lobaria pulmonaria"
$ws.Range("H23").Value = "(L.) Hoffm."
$ws.Range("Q23").Value = 547783
$ws.Range("R23").Value = 6960015
$ws.Range("Z23").ClearContents()
$ws.Range("AB23").ClearContents()
$ws.Range("AC23").ClearContents()

# Row 24
$ws.Range("A24").Value = 112043806
$ws.Range("B24").Value = 56398
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("Q24").Value = 547716
$ws.Range("R24").Value = 6960073
$ws.Range("Z24").ClearContents()
$ws.Range("AB24").ClearContents()
$ws.Range("AC24").Value = "ringhack äldre"

# Row 25
$ws.Range("A25").Value = 112043836
$ws.Range("B25").Value = 78578
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 6458
$ws.Range("F25").Value = "Lunglav"
$ws.Range("G25").Value = "Lobaria pulmonaria"
$ws.Range("H25").Value = "(L.) Hoffm."
$ws.Range("Q25").Value = 547781
$ws.Range("R25").Value = 6960012
$ws.Range("Z25").ClearContents()
$ws.Range("AB25").ClearContents()
$ws.Range("AC25").ClearContents()

# Row 26
$ws.Range("A26").Value = 112043828
$ws.Range("B26").Value = 78578
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6458
$ws.Range("F26").Value = "Lunglav"
$ws.Range("G26").Value = "Lobaria pulmonaria"
$ws.Range("H26").Value = "(L.) Hoffm."
$ws.Range("Q26").Value = 547484
$ws.Range("R26").Value = 6960144
$ws.Range("Z26").ClearContents()
$ws.Range("AB26").ClearContents()
$ws.Range("AC26").ClearContents()

# Row 27
$ws.Range("A27").Value = 112043851
$ws.Range("B27").Value = 89965
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 760
$ws.Range("F27").Value = "Doftticka"
$ws.Range("G27").Value = "Haploporus odorus"
$ws.Range("H27").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q27").Value = 547815
$ws.Range("R27").Value = 6960221
$ws.Range("Z27").ClearContents()
$ws.Range("AB27").ClearContents()
$ws.Range("AC27").ClearContents()

# Row 28
$ws.Range("A28").Value = 112043811
$ws.Range("B28").Value = 90332
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 4769
$ws.Range("F28").Value = "Svavelriska"
$ws.Range("G28").Value = "Lactarius scrobiculatus"
$ws.Range("H28").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q28").Value = 547793
$ws.Range("R28").Value = 6960088
$ws.Range("Z28").ClearContents()
$ws.Range("AB28").ClearContents()
$ws.Range("AC28").ClearContents()

# Row 29
$ws.Range("A29").Value = 112043853
$ws.Range("B29").Value = 77515
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("Q29").Value = 547529
$ws.Range("R29").Value = 6960183
$ws.Range("Z29").ClearContents()
$ws.Range("AB29").ClearContents()
$ws.Range("K29").ClearContents()
$ws.Range("L29").ClearContents()
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("AC29").ClearContents()

# Row 30
$ws.Range("A30").Value = 112043823
$ws.Range("B30").Value = 78578
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6458
$ws.Range("F30").Value = "Lunglav"
$ws.Range("G30").Value = "Lobaria pulmonaria"
$ws.Range("H30").Value = "(L.) Hoffm."
$ws.Range("Q30").Value = 547553
$ws.Range("R30").Value = 6960101
$ws.Range("Z30").ClearContents()
$ws.Range("AB30").ClearContents()
$ws.Range("AC30").ClearContents()

# Row 31
$ws.Range("A31").Value = 112043832
$ws.Range("B31").Value = 78578
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 6458
$ws.Range("F31").Value = "Lunglav"
$ws.Range("G31").Value = "Lobaria pulmonaria"
$ws.Range("H31").Value = "(L.) Hoffm."
$ws.Range("Q31").Value = 547489
$ws.Range("R31").Value = 6960079
$ws.Range("Z31").ClearContents()
$ws.Range("AB31").ClearContents()
$ws.Range("AC31").ClearContents()

# Row 32
$ws.Range("A32").Value = 112043818
$ws.Range("B32").Value = 78578
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 6458
$ws.Range("F32").Value = "Lunglav"
$ws.Range("G32").Value = "Lobaria pulmonaria"
$ws.Range("H32").Value = "(L.) Hoffm."
$ws.Range("Q32").Value = 547774
$ws.Range("R32").Value = 6960191
$ws.Range("Z32").ClearContents()
$ws.Range("AB32").ClearContents()
$ws.Range("AC32").ClearContents()

# Row 33
$ws.Range("A33").Value = 112043800
$ws.Range("B33").Value = 86223
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 4412
$ws.Range("F33").Value = "Äggvaxskivling"
$ws.Range("G33").Value = "Hygrophorus karstenii"
$ws.Range("H33").Value = "Sacc. & Cub."
$ws.Range("Q33").Value = 547427
$ws.Range("R33").Value = 6960212
$ws.Range("Z33").ClearContents()
$ws.Range("AB33").ClearContents()
$ws.Range("AC33").ClearContents()

# Row 34
$ws.Range("A34").Value = 112043843
$ws.Range("B34").Value = 90687
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 5964
$ws.Range("F34").Value = "Fjällig taggsvamp s.str."
$ws.Range("G34").Value = "Sarcodon imbricatus s.str."
$ws.Range("H34").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q34").Value = 547462
$ws.Range("R34").Value = 6960196
$ws.Range("Z34").ClearContents()
$ws.Range("AB34").ClearContents()
$ws.Range("AC34").ClearContents()

# Row 35
$ws.Range("A35").Value = 112043842
$ws.Range("B35").Value = 90687
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 5964
$ws.Range("F35").Value = "Fjällig taggsvamp s.str."
$ws.Range("G35").Value = "Sarcodon imbricatus s.str."
$ws.Range("H35").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q35").Value = 547638
$ws.Range("R35").Value = 6960094
$ws.Range("Z35").ClearContents()
$ws.Range("AB35").ClearContents()
$ws.Range("AC35").ClearContents()

# Row 36
$ws.Range("A36").Value = 112043847
$ws.Range("B36").Value = 90687
$ws.Range("D36").Value = "LC"
$ws.Range("E36").Value = 5964
$ws.Range("F36").Value = "Fjällig taggsvamp s.str."
$ws.Range("G36").Value = "Sarcodon imbricatus s.str."
$ws.Range("H36").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q36").Value = 547492
$ws.Range("R36").Value = 6960149
$ws.Range("Z36").ClearContents()
$ws.Range("AB36").ClearContents()
$ws.Range("AC36").ClearContents()

# Row 37
$ws.Range("A37").Value = 112043841
$ws.Range("B37").Value = 90687
$ws.Range("D37").Value = "LC"
$ws.Range("E37").Value = 5964
$ws.Range("F37").Value = "Fjällig taggsvamp s.str."
$ws.Range("G37").Value = "Sarcodon imbricatus s.str."
$ws.Range("H37").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q37").Value = 547829
$ws.Range("R37").Value = 6960095
$ws.Range("Z37").ClearContents()
$ws.Range("AB37").ClearContents()
$ws.Range("AC37").ClearContents()

# Row 38
$ws.Range("A38").Value = 112043840
$ws.Range("B38").Value = 90687
$ws.Range("D38").Value = "LC"
$ws.Range("E38").Value = 5964
$ws.Range("F38").Value = "Fjällig taggsvamp s.str."
$ws.Range("G38").Value = "Sarcodon imbricatus s.str."
$ws.Range("H38").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q38").Value = 547722
$ws.Range("R38").Value = 6960170
$ws.Range("Z38").ClearContents()
$ws.Range("AB38").ClearContents()
$ws.Range("AC38").ClearContents()

# Row 39
$ws.Range("A39").Value = 112043845
$ws.Range("B39").Value = 90687
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 5964
$ws.Range("F39").Value = "Fjällig taggsvamp s.str."
$ws.Range("G39").Value = "Sarcodon imbricatus s.str."
$ws.Range("H39").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q39").Value = 547639
$ws.Range("R39").Value = 6960116
$ws.Range("Z39").ClearContents()
$ws.Range("AB39").ClearContents()
$ws.Range("AC39").ClearContents()

# Row 40
$ws.Range("A40").Value = 112043848
$ws.Range("B40").Value = 90687
$ws.Range("D40").Value = "LC"
$ws.Range("E40").Value = 5964
$ws.Range("F40").Value = "Fjällig taggsvamp s.str."
$ws.Range("G40").Value = "Sarcodon imbricatus s.str."
$ws.Range("H40").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q40").Value = 547485
$ws.Range("R40").Value = 6960143
$ws.Range("Z40").ClearContents()
$ws.Range("AB40").ClearContents()
$ws.Range("AC40").ClearContents()
